# "Generate Report for Handback" — localization-status.xlsx
#
# The handback process completed for the d1c77145-... source file in both
# target locales (zh-cn, de-de): the status moves from "Ready for handoff"
# to "Handed back: in sync with en-US", the generated target/handback files
# are recorded, and the handback timestamp is stamped. The Overview sheet
# mirrors the per-locale status. Columns that now hold longer strings are
# widened to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile    = "d1c77145-b09c-44b4-80e2-68135d36d963.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7b743326eac3d380f4661dbcb30033c699971ed/e2e/d1c77145-b09c-44b4-80e2-68135d36d963.md"

# ---------------------------------------------------------------------
# Overview sheet: reflect the new status for both locale columns (E, F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = $newStatus
$wsOverview.Cells.Item(2, 6).Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status
$wsZhCn.Cells.Item(2, 3).Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14

# Latest Target File (I2) -> hyperlink to the handed-back source doc
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2, 9), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)

# Latest Handback File (J2)
$wsZhCn.Cells.Item(2, 10).Value = "d1c77145-b09c-44b4-80e2-68135d36d963.4c5163c3fae4d8e30d64e079ba3c0f7c5a6cb013.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZhCn.Cells.Item(2, 11).Value = "2016-08-16 20:54:59"

$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status
$wsDeDe.Cells.Item(2, 3).Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 29.14

# Latest Target File (I2) -> hyperlink to the handed-back source doc
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2, 9), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)

# Latest Handback File (J2)
$wsDeDe.Cells.Item(2, 10).Value = "d1c77145-b09c-44b4-80e2-68135d36d963.4c5163c3fae4d8e30d64e079ba3c0f7c5a6cb013.de-de.xlf"

# Latest Handback DateTime (K2)
$wsDeDe.Cells.Item(2, 11).Value = "2016-08-16 20:55:15"

$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
